# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP) for rows 2-16.
# Column F (Win) is unchanged. Column G (sum) = B + C + D + E.
$data = @{
    2  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    3  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    4  = @(3.286832544864788, 1.655778082260271, 3.537761648806719,  0.4942365360607697)
    5  = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697)
    6  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
    7  = @(3.286832544864788, 10.34677158129881, 0.1494219747398047, 10.19245300693656)
    8  = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 10.19245300693656)
    9  = @(1.455362044514542, 1.655778082260271, 3.537761648806719,  0.4942365360607697)
    10 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    11 = @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697)
    12 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697)
    13 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    14 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    15 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
    16 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = ($b + $c + $d + $e)
}
